$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.926.63"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "3.104.89"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'576.25"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Value = "'177.66"
$ws.Range("E6").Value = "  +2.90%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "3.104.11"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("E10").Value = "  -2.31%  "

$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").Value = "'0.468"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("D14").Value = "'36.14"
$ws.Range("E14").Value = "  -1.28%  "

$ws.Range("D15").Value = "'0.121"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").Value = "3.621.70"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").Value = "66.918.07"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "'7.04"

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'16.76"
$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.105.47"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "'480.02"
$ws.Range("E21").Value = "  -1.93%  "

$ws.Range("D22").Value = "'7.79"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").Value = "'0.690"
$ws.Range("E23").Value = "  -1.34%  "

$ws.Range("D24").Value = "'83.68"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").Value = "'12.56"
$ws.Range("E25").Value = "  -3.41%  "

$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  -1.86%  "

$ws.Range("D27").Value = "'10.08"
$ws.Range("E27").Value = "  -4.10%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").Value = "'7.90"
$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("E30").Value = "  -1.96%  "

$ws.Range("D31").Value = "'2.60"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("D32").Value = "'28.02"
$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("D34").Value = "0.0₃0942"
$ws.Range("E34").Value = "  +1.85%  "

$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("E36").Value = "  +3.62%  "

$ws.Range("D37").Value = "'5.59"
$ws.Range("E37").Value = "  -3.27%  "

$ws.Range("D38").Value = "'0.944"
$ws.Range("E38").Value = "  -3.04%  "

$ws.Range("D39").Value = "'0.311"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("D40").Value = "'49.05"
$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").Value = "'8.33"
$ws.Range("E43").Value = "  -1.27%  "

$ws.Range("D44").Value = "'2.70"
$ws.Range("E44").Value = "  +5.51%  "

$ws.Range("D45").Value = "2.799.35"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "'371.44"
$ws.Range("E46").Value = "  -3.32%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'135.65"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0344"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "  +2.11%  "
